{"js": "// Updated the User Stories\n//\n// The paragraph that begins \"La pareja formada por Jos\u00e9 Manuel y Jos\u00e9...\"\n// (followed by the merged \"Con respecto a las pruebas...\" sentence) is\n// rewritten as a single run of text, gains a first-line indent matching\n// the rest of the document's body paragraphs, and keeps its (empty)\n// \"_GoBack\" bookmark, which now sits at the very start of the paragraph\n// content instead of in the middle of the old run sequence.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the target paragraph by its distinctive leading text instead of\n// a hard-coded index, so the script is resilient to minor structural\n// differences.\nconst marker = \"La pareja formada por Jos\u00e9 Manuel y Jos\u00e9 tuvo\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(marker) === 0) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the 'La pareja formada...' paragraph.\");\n}\n\ntarget.load(\"text\");\nawait context.sync();\nconst fullText = target.text;\n\n// 1) Give the paragraph the same first-line indent (708 twips = 35.4pt)\n//    used by the other body paragraphs in this document.\ntarget.firstLineIndent = 35.4;\n\n// 2) Collapse all of the paragraph's runs (and the old mid-paragraph\n//    \"_GoBack\" bookmark) into a single run containing the same text.\nconst paraRange = target.getRange();\nparaRange.insertText(fullText, Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) Re-create the \"_GoBack\" bookmark at the start of the paragraph.\nconst freshParagraphs = body.paragraphs;\nfreshParagraphs.load(\"items/text\");\nawait context.sync();\n\nlet refreshedTarget = null;\nfor (let i = 0; i < freshParagraphs.items.length; i++) {\n  if (freshParagraphs.items[i].text.indexOf(marker) === 0) {\n    refreshedTarget = freshParagraphs.items[i];\n    break;\n  }\n}\n\nconst startRange = refreshedTarget.getRange(Word.RangeLocation.start);\nstartRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Updated the User Stories\n#\n# The paragraph that begins \"La pareja formada por Jos\u00e9 Manuel y Jos\u00e9...\"\n# (followed by the merged \"Con respecto a las pruebas...\" sentence) is\n# rewritten as a single run of text, gains a first-line indent matching\n# the rest of the document's body paragraphs, and keeps its (empty)\n# \"_GoBack\" bookmark, which now sits at the very start of the paragraph\n# content instead of in the middle of the old run sequence.\n\n$d = $word.ActiveDocument\n\n# Locate the target paragraph by its distinctive leading text instead of\n# a hard-coded index, so the script is resilient to minor structural\n# differences.\n$marker = \"La pareja formada por\"\n$targetIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.StartsWith($marker)) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not locate the 'La pareja formada...' paragraph.\"\n}\n\n$p = $d.Paragraphs.Item($targetIndex)\n\n# 1) Give the paragraph the same first-line indent (708 twips = 35.4pt)\n#    used by the other body paragraphs in this document.\n$p.Format.FirstLineIndent = 35.4\n\n# Capture the full paragraph text (without the trailing paragraph mark)\n# before collapsing the runs / removing the mid-paragraph bookmark.\n$fullRange = $p.Range\n$fullRange.MoveEnd(1, -1) | Out-Null   # wdCharacter = 1; drop the pilcrow\n$fullText = $fullRange.Text\n\n# 2) Collapse all of the paragraph's runs (and the old mid-paragraph\n#    \"_GoBack\" bookmark) down to nothing, then insert the same text back\n#    as a single run.\n$fullRange.Delete()\n\n$p2 = $d.Paragraphs.Item($targetIndex)\n$insertRange = $p2.Range\n$insertRange.MoveEnd(1, -1) | Out-Null\n$insertRange.InsertAfter($fullText)\n\n# 3) Re-create the (empty/collapsed) \"_GoBack\" bookmark at the very start\n#    of the paragraph.\n$p3 = $d.Paragraphs.Item($targetIndex)\n$startRange = $p3.Range\n$startRange.MoveEnd(1, -1) | Out-Null\n$startRange.Collapse(1)   # wdCollapseStart = 1\n$d.Bookmarks.Add(\"_GoBack\", $startRange) | Out-Null\n"}
